$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26, shifting existing rows 26-30 down to 27-31
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with data (copy pattern values + new specifics from the diff)
$ws.Cells.Item(26, 1).Value = 11
$ws.Cells.Item(26, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value = "Bíobío"
$ws.Cells.Item(26, 4).Value = 44782
$ws.Cells.Item(26, 4).NumberFormat = $ws.Cells.Item(27, 4).NumberFormat
$ws.Cells.Item(26, 5).Value = 8
$ws.Cells.Item(26, 6).Value = 100112026
$ws.Cells.Item(26, 7).Value = "Haba"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 30
$ws.Cells.Item(26, 11).Value = 15000
$ws.Cells.Item(26, 12).Value = 15000
$ws.Cells.Item(26, 13).Value = 15000
$ws.Cells.Item(26, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(26, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(26, 16).Value = 600
$ws.Cells.Item(26, 17).Value = 25
$ws.Cells.Item(26, 18).Value = "Hortaliza"
